$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: Wins / Losses / Ties in columns AD, AE, AF (row 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, centered, bordered) from an existing header cell (AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-43: team record values
$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 102
    $ws.Cells.Item($r, 31).Value = 60
    $ws.Cells.Item($r, 32).Value = 0
}
